$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 rework: the old "SAN DIEGO AREA TOTALS" label lived in B2 (Calibri 12,
# same look-and-feel as the other section headers). It now moves one column to the
# left into A2, and B2 picks up a plain "Totals" label in the sheet's default font. ---

# 1) Carry B2's current (Calibri 12 / style index 1) formatting over to A2 first,
#    while B2 still has it, via copy / paste-special (formats only).
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 2) Now write the actual values: A2 gets the label that used to sit in B2, and
#    B2 becomes "Totals" (reusing the shared string already used elsewhere, e.g. B18/B25/B28).
$ws.Range("A2").Value2 = "SAN DIEGO AREA TOTALS "
$ws.Range("B2").Value2 = "Totals"

# 3) B2 should end up with the sheet's default formatting (no explicit style), so
#    stamp it with a throwaway default-formatted cell's formats, then clean that helper up.
$ws.Range("F1").Value2 = "tmp"
$ws.Range("F1").Copy()
$ws.Range("B2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("F1").Clear()

# --- Column widths: column A now holds the same long header text that column B
# used to, so widen it to line up with column B. ---
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# --- Selection moves from D15 to A6. ---
$ws.Range("A6").Select()
